# Sort the "NEW YORK" listings sheet by "Price per Sqft" (column E), ascending,
# keeping the header row in place (row 1) and reordering rows 2:36 together.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NEW YORK")
$ws.Activate()

# Full table including the header row.
$sortRange = $ws.Range("A1:I36")

# Sort key = "Price per Sqft" column.
$keyRange = $ws.Range("E2:E36")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange, 0, 1, 0, 0)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 1
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1
$ws.Sort.Apply()
